# The workbook originally has two sheets, in this order:
#   1) hotel_info  - one header row + one data row describing a hotel
#   2) review_info - one header row (no data rows) describing reviews
#
# This edit:
#   1. Inserts a new "State" column into hotel_info, between "Hotel_Name"
#      and "City", with the value "Louisiana" for the existing hotel row.
#   2. Swaps the tab order so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# Insert a new column C ("State") in hotel_info, shifting City (and
# everything after it) one column to the right.
$wsHotel.Range("C:C").Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# Move hotel_info to sit right after review_info, i.e. swap their order.
$wsHotel.Move($null, $wsReview)
